$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.825.39"
$ws.Range("D3").Value = "1.628.28"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.26"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.20"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0879"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.859.57"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "1.620.24"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.86"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").Value = "27.844.24"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.19"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.34"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.94"
$ws.Range("E23").Value = "  -5.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.29"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.47"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").Value = "1.411.97"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.554"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.65"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "1.768.81"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.47"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.61"
$ws.Range("E51").Value = "  +0.38%  "
